$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.445.16'
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').Value = '2.090.88'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '228.21'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.01%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.613'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.76%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '60.89'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  -0.01%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.382'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.37%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0837'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').Value = '2.400.55'
$ws.Range('E12').Value = '  +2.36%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '14.85'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.28%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '22.39'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +6.33%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.786'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.08%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.44'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +4.23%  '
$ws.Range('D17').Value = '2.093.37'
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').Value = '38.375.83'
$ws.Range('E18').Value = '  +1.50%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '71.28'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +2.45%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.06'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('D21').Value = '0.0₃0834'
$ws.Range('E21').Value = '  +1.24%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '225.60'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('E23').Value = '  -0.03%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('E25').Value = '  +2.10%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '170.02'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.15%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.43'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.96%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.135'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +4.71%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '19.06'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.31%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.37'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +8.68%  '
$ws.Range('E31').Value = '  -0.48%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.33'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +5.18%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.81'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +6.45%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.51'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +2.76%  '
$ws.Range('E35').Value = '  +0.92%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '6.45'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -2.38%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +2.06%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.56'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +2.44%  '
$ws.Range('E39').Value = '  +0.00%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '18.47'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('D41').Value = '1.540.35'
$ws.Range('E41').Value = '  -0.01%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '100.05'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.86%  '
$ws.Range('E43').Value = '  +1.44%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0930'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('E45').Value = '  -0.25%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '7.78'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +10.10%  '
$ws.Range('E47').Value = '  -2.08%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.12'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').Value = '2.287.17'
$ws.Range('E51').Value = '  +2.35%  '
